# Estado de cuenta update:
#  - updates Valor Mora total and Cant. Periodos counters
#  - re-sorts the "Periodo Mora" rows so period "2311" (with its special
#    partial value) moves to the top of the list and the remaining periods
#    follow in ascending calendar order
#  - appends a new period row ("2508") at the end of the workers table,
#    shifting the trailer ("firma") rows down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header figures -------------------------------------------------
$ws.Range("E11").Value = 1039147   # VALOR MORA (total)
$ws.Range("F13").Value = 23        # Cant. Periodos

# --- make room for the new period row --------------------------------
# Row 37 is currently the last data row (styled with the heavier bottom
# border). Insert a new row after it so the table grows by one entry;
# this naturally pushes the two signature/trailer rows down by one row.
$ws.Rows("38:38").Insert()

# New row 38 should look like the old "last row" (row 37 before the
# insert); copy that formatting down into it.
$ws.Range("B37:J37").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)   # xlPasteFormats

# Row 37 is no longer the last row of the table, so it should now match
# the regular interior row formatting (same as row 36).
$ws.Range("B36:J36").Copy()
$ws.Range("B37:J37").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- fill the period rows in ascending order -------------------------
$periods = @("2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")

$r = 17
foreach ($p in $periods) {
    $ws.Range("B" + $r).Value = "CC"
    $ws.Range("C" + $r).Value = "1007621086"
    $ws.Range("D" + $r).Value = "LEONARDO FAVIO VILLEGAS TORRES"
    $ws.Range("E" + $r).Value = $p
    if ($p -eq "2311") {
        $ws.Range("F" + $r).Value = 24747
    } else {
        $ws.Range("F" + $r).Value = 46400
    }
    $ws.Range("G" + $r).Value = 1160000
    $r = $r + 1
}
